# ThreadObjekte.xlsx re-edit:
#   - move the active selection to B19 (was B31)
#   - widen column A to fit the (longer) "Grundlagen"/"Anforderungen" labels
#   - set the print page setup to A4 portrait (paperSize=9, portrait orientation)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (was 15.5703125 chars, now considerably wider to fit longer text)
$ws.Columns.Item(1).ColumnWidth = 24.5

# Move / update the current selection to B19
$ws.Range("B19").Select()

# Page setup: A4, portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
